$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws.Range("B2").Value = 0.6336537368188535
$ws.Range("B3").Value = 41.02581503960927
$ws.Range("B4").Value = 42.39362216031643
$ws.Range("B5").Value = 43.05267931637545
$ws.Range("B6").Value = 44.07520410564044
$ws.Range("B7").Value = 43.96215999009672
$ws.Range("B8").Value = 42.14904433939637
$ws.Range("B9").Value = 41.27835505142821
$ws.Range("B10").Value = 41.88006630321046
$ws.Range("B11").Value = 45.4139720420827
$ws.Range("B12").Value = 49.00096300291401
$ws.Range("B13").Value = 61.13025967955659
$ws.Range("B14").Value = 61.18336750180829
$ws.Range("B15").Value = 68.99110855652924
$ws.Range("B16").Value = 73.97153386361917
$ws.Range("B17").Value = 57.52799620924602
$ws.Range("B18").Value = 58.75507011031424
$ws.Range("B19").Value = 60.95253353026607
$ws.Range("B20").Value = 61.23936655450935
$ws.Range("B21").Value = 61.75224678495422
$ws.Range("B22").Value = 60.83731037154546
$ws.Range("B23").Value = 58.88551498428396
$ws.Range("B24").Value = 58.34734598421247
$ws.Range("B25").Value = 57.98179921010325
$ws.Range("B26").Value = 56.99148187264407
$ws.Range("B27").Value = 55.55412484472498
$ws.Range("B28").Value = 56.70954534706966
$ws.Range("B29").Value = 57.77178093952842
$ws.Range("B30").Value = 58.35988993994804
$ws.Range("B31").Value = 57.43851550461792
$ws.Range("B32").Value = 59.77871294774496
$ws.Range("B33").Value = 58.24657571219976
$ws.Range("B34").Value = 58.2018230531539
$ws.Range("B35").Value = 59.36725624646149
$ws.Range("B36").Value = 60.31241172270381
$ws.Range("B37").Value = 61.88462775052184
$ws.Range("B38").Value = 62.94152091356012

$ws = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws.Range("B2").Value = 63.837349565088
$ws.Range("C2").Value = 56.14261451516851
$ws.Range("D2").Value = 71.53208461500748
$ws.Range("B3").Value = 64.47100330190685
$ws.Range("C3").Value = 53.5890046354431
$ws.Range("D3").Value = 75.35300196837059
$ws.Range("B4").Value = 65.1046570387257
$ws.Range("C4").Value = 51.77698498148412
$ws.Range("D4").Value = 78.43232909596729
$ws.Range("B5").Value = 65.73831077554456
$ws.Range("C5").Value = 50.34884067570559
$ws.Range("D5").Value = 81.12778087538354
$ws.Range("B6").Value = 66.37196451236342
$ws.Range("C6").Value = 49.16601387189321
$ws.Range("D6").Value = 83.57791515283363

$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Range("B2").Value = 0.4631494590506633
$ws.Range("B3").Value = 40.85531076184108
$ws.Range("B4").Value = 42.22311788254824
$ws.Range("B5").Value = 42.88217503860727
$ws.Range("B6").Value = 43.90469982787225
$ws.Range("B7").Value = 43.79165571232853
$ws.Range("B8").Value = 41.97854006162818
$ws.Range("B9").Value = 41.10785077366002
$ws.Range("B10").Value = 41.70956202544227
$ws.Range("B11").Value = 45.24346776431451
$ws.Range("B12").Value = 48.83045872514582
$ws.Range("B13").Value = 60.95975540178841
$ws.Range("B14").Value = 61.0128632240401
$ws.Range("B15").Value = 68.82060427876105
$ws.Range("B16").Value = 73.80102958585098
$ws.Range("B17").Value = 57.35749193147783
$ws.Range("B18").Value = 58.58456583254605
$ws.Range("B19").Value = 60.78202925249789
$ws.Range("B20").Value = 61.06886227674116
$ws.Range("B21").Value = 61.58174250718603
$ws.Range("B22").Value = 60.66680609377728
$ws.Range("B23").Value = 58.71501070651577
$ws.Range("B24").Value = 58.17684170644428
$ws.Range("B25").Value = 57.81129493233506
$ws.Range("B26").Value = 56.82097759487588
$ws.Range("B27").Value = 55.38362056695679
$ws.Range("B28").Value = 56.53904106930147
$ws.Range("B29").Value = 57.60127666176023
$ws.Range("B30").Value = 58.18938566217985
$ws.Range("B31").Value = 57.26801122684973
$ws.Range("B32").Value = 59.60820866997677
$ws.Range("B33").Value = 58.07607143443157
$ws.Range("B34").Value = 58.03131877538571
$ws.Range("B35").Value = 59.1967519686933
$ws.Range("B36").Value = 60.14190744493562
$ws.Range("B37").Value = 61.71412347275365
$ws.Range("B38").Value = 62.77101663579193
$ws.Range("B39").Value = 63.66684528731981
$ws.Range("B40").Value = 63.80182231980051
$ws.Range("B41").Value = 63.11649980607323
$ws.Range("B42").Value = 63.53099871295088
$ws.Range("B43").Value = 63.25836073300315

$ws = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws.Range("B2").Value = 59.84443858291828
$ws.Range("C2").Value = 52.517815811601
$ws.Range("D2").Value = 67.17106135423555
$ws.Range("B3").Value = 60.30758804196894
$ws.Range("C3").Value = 49.9461787523805
$ws.Range("D3").Value = 70.66899733155738
$ws.Range("B4").Value = 60.7707375010196
$ws.Range("C4").Value = 48.080654613207
$ws.Range("D4").Value = 73.46082038883222
$ws.Range("B5").Value = 61.23388696007027
$ws.Range("C5").Value = 46.58064141743571
$ws.Range("D5").Value = 75.88713250270482
$ws.Range("B6").Value = 61.69703641912093
$ws.Range("C6").Value = 45.3142098569576
$ws.Range("D6").Value = 78.07986298128425
